$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1737.6522
$ws.Range("I100").Value = 1328.8889
$ws.Range("J100").Value = 2000.4286
$ws.Range("K100").Value = 1328.8889
$ws.Range("L100").Value = 2000.4286
$ws.Range("M100").Value = -787.8888999999999
$ws.Range("N100").Value = -3082.4286
$ws.Range("H129").Value = 839.9231
$ws.Range("I129").Value = 307.45456
$ws.Range("J129").Value = 1049.1072
$ws.Range("K129").Value = 922.36368
$ws.Range("L129").Value = 3147.3216
$ws.Range("M129").Value = 4077.63632
$ws.Range("N129").Value = -13147.3216
$ws.Range("H137").Value = 1097.5574
$ws.Range("I137").Value = 844.5161000000001
$ws.Range("J137").Value = 1359.0333
$ws.Range("K137").Value = 2533.5483
$ws.Range("L137").Value = 4077.0999
$ws.Range("M137").Value = 16.45169999999962
$ws.Range("N137").Value = -9177.099900000001
$ws.Range("H138").Value = 1361.4
$ws.Range("I138").Value = 872.7826
$ws.Range("J138").Value = 1507.3507
$ws.Range("K138").Value = 2618.3478
$ws.Range("L138").Value = 4522.0521
$ws.Range("M138").Value = 2521.6522
$ws.Range("N138").Value = -14802.0521

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3335.1094
$ws.Range("I32").Value = 3079.7903
$ws.Range("K32").Value = 3079.7903
$ws.Range("M32").Value = -2792.7903
$ws.Range("H45").Value = 1462.4
$ws.Range("I45").Value = 1446.2858
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1446.2858
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -1069.2858
$ws.Range("N45").Value = -2254
$ws.Range("H74").Value = 1517.7778
$ws.Range("I74").Value = 955.38464
$ws.Range("J74").Value = 2980
$ws.Range("K74").Value = 955.38464
$ws.Range("L74").Value = 2980
$ws.Range("M74").Value = -81.38463999999999
$ws.Range("N74").Value = -4728
$ws.Range("H77").Value = 1517.7778
$ws.Range("I77").Value = 955.38464
$ws.Range("J77").Value = 2980
$ws.Range("K77").Value = 4776.9232
$ws.Range("L77").Value = 14900
$ws.Range("M77").Value = -408.9232000000002
$ws.Range("N77").Value = -23636
$ws.Range("H132").Value = 1404.0785
$ws.Range("I132").Value = 1136.775
$ws.Range("J132").Value = 2376.0908
$ws.Range("K132").Value = 3410.325
$ws.Range("L132").Value = 7128.2724
$ws.Range("M132").Value = -880.3250000000003
$ws.Range("N132").Value = -12188.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 65000
$ws.Range("J62").Value = 65000
$ws.Range("L62").Value = 65000
$ws.Range("N62").Value = -66372
$ws.Range("H65").Value = 65000
$ws.Range("J65").Value = 65000
$ws.Range("L65").Value = 195000
$ws.Range("N65").Value = -201864
$ws.Range("H99").Value = 66667730
$ws.Range("I99").Value = 76924100
$ws.Range("K99").Value = 76924100
$ws.Range("M99").Value = -76922602
$ws.Range("H128").Value = 1040
$ws.Range("I128").Value = 1040
$ws.Range("K128").Value = 3120
$ws.Range("M128").Value = -630

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2756.5715
$ws.Range("I31").Value = 2891.6924
$ws.Range("K31").Value = 2891.6924
$ws.Range("M31").Value = -2596.6924
$ws.Range("H34").Value = 2756.5715
$ws.Range("I34").Value = 2891.6924
$ws.Range("K34").Value = 2891.6924
$ws.Range("M34").Value = -2689.6924
$ws.Range("H58").Value = 643.0943600000001
$ws.Range("I58").Value = 586.3946999999999
$ws.Range("J58").Value = 786.73334
$ws.Range("K58").Value = 586.3946999999999
$ws.Range("L58").Value = 786.73334
$ws.Range("M58").Value = -383.3946999999999
$ws.Range("N58").Value = -1192.73334
$ws.Range("H107").Value = 1472.2858
$ws.Range("I107").Value = 2096.2856
$ws.Range("J107").Value = 848.2857
$ws.Range("K107").Value = 2096.2856
$ws.Range("L107").Value = 848.2857
$ws.Range("M107").Value = -176.2856000000002
$ws.Range("N107").Value = -4688.2857
$ws.Range("H132").Value = 3114.3809
$ws.Range("I132").Value = 2712.75
$ws.Range("K132").Value = 8138.25
$ws.Range("M132").Value = -5608.25
$ws.Range("H134").Value = 1038.6857
$ws.Range("I134").Value = 1059.1428
$ws.Range("K134").Value = 3177.4284
$ws.Range("M134").Value = -642.4284000000002
$ws.Range("H136").Value = 643.0943600000001
$ws.Range("I136").Value = 586.3946999999999
$ws.Range("J136").Value = 786.73334
$ws.Range("K136").Value = 1759.1841
$ws.Range("L136").Value = 2360.20002
$ws.Range("M136").Value = 790.8159000000001
$ws.Range("N136").Value = -7460.20002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9891.412
$ws.Range("J3").Value = 15740.5
$ws.Range("L3").Value = 47221.5
$ws.Range("N3").Value = -47445.5
$ws.Range("H4").Value = 559320.9399999999
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 450
$ws.Range("M4").Value = -338
$ws.Range("H113").Value = 698.5806
$ws.Range("I113").Value = 599.5
$ws.Range("K113").Value = 1798.5
$ws.Range("M113").Value = 371.5
$ws.Range("H121").Value = 787.8333
$ws.Range("J121").Value = 999.25
$ws.Range("L121").Value = 2997.75
$ws.Range("N121").Value = -5617.75
$ws.Range("H131").Value = 27028326
$ws.Range("I131").Value = 76923300
$ws.Range("J131").Value = 1884.5416
$ws.Range("K131").Value = 230769900
$ws.Range("L131").Value = 5653.6248
$ws.Range("M131").Value = -230764860
$ws.Range("N131").Value = -15733.6248
$ws.Range("H138").Value = 2169.0527
$ws.Range("I138").Value = 2110.5833
$ws.Range("J138").Value = 2269.2856
$ws.Range("K138").Value = 6331.749899999999
$ws.Range("L138").Value = 6807.8568
$ws.Range("M138").Value = -1191.749899999999
$ws.Range("N138").Value = -17087.8568
$ws.Range("H139").Value = 1906.2703
$ws.Range("I139").Value = 2116.7368
$ws.Range("K139").Value = 6350.2104
$ws.Range("M139").Value = -1210.2104
$ws.Range("H140").Value = 23100.314
$ws.Range("I140").Value = 54882.76
$ws.Range("J140").Value = 2875.121
$ws.Range("K140").Value = 164648.28
$ws.Range("L140").Value = 8625.363000000001
$ws.Range("M140").Value = -159468.28
$ws.Range("N140").Value = -18985.363
$ws.Range("H141").Value = 2590.182
$ws.Range("I141").Value = 2645.9
$ws.Range("K141").Value = 7937.700000000001
$ws.Range("M141").Value = -2757.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2187.8462
$ws.Range("I132").Value = 1601.8889
$ws.Range("K132").Value = 4805.6667
$ws.Range("M132").Value = -2275.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1750
$ws.Range("J46").Value = 2125
$ws.Range("L46").Value = 2125
$ws.Range("N46").Value = -2501
$ws.Range("H122").Value = 16669038
$ws.Range("I122").Value = 22729530
$ws.Range("J122").Value = 2685
$ws.Range("K122").Value = 68188590
$ws.Range("L122").Value = 8055
$ws.Range("M122").Value = -68186140
$ws.Range("N122").Value = -12955
$ws.Range("H132").Value = 23250.191
$ws.Range("I132").Value = 1531.2916
$ws.Range("J132").Value = 45913.39
$ws.Range("K132").Value = 4593.8748
$ws.Range("L132").Value = 137740.17
$ws.Range("M132").Value = -2063.8748
$ws.Range("N132").Value = -142800.17
$ws.Range("H136").Value = 1368.4546
$ws.Range("I136").Value = 1132.4117
$ws.Range("J136").Value = 2171
$ws.Range("K136").Value = 3397.2351
$ws.Range("L136").Value = 6513
$ws.Range("M136").Value = -847.2351000000003
$ws.Range("N136").Value = -11613

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 35716460
$ws.Range("J122").Value = 1900
$ws.Range("L122").Value = 5700
$ws.Range("N122").Value = -10600
